$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: new section title (mirrors row 6/13 style "Scenario ..." label in col K) ---
$ws.Range("K21").Value = "Scenario E, 5% sellers, TD=3"

# --- Row 22: column headers for % Eff / % Liquidity / No Seller (copy header style from row 7) ---
$ws.Range("K7:M7").Copy()
$ws.Range("K22:M22").PasteSpecial(-4122)
$ws.Range("K22").Value = "% Eff"
$ws.Range("L22").Value = "% Liquidity"
$ws.Range("M22").Value = "No Seller"

# --- Rows 23-27: data rows (copy number formats from an existing data row, then set values) ---
$ws.Range("K8:N8").Copy()
$ws.Range("K23:N23").PasteSpecial(-4122)
$ws.Range("K24:N24").PasteSpecial(-4122)
$ws.Range("K25:N25").PasteSpecial(-4122)
$ws.Range("K26:N26").PasteSpecial(-4122)
$ws.Range("K27:N27").PasteSpecial(-4122)

$ws.Range("J23").Value = " "
$ws.Range("K23").Value = 0.108
$ws.Range("L23").Value = 0.53680000000000005
$ws.Range("M23").Value = 0.3553
$ws.Range("N23").Formula = "=SUM(K23:M23)"

$ws.Range("J24").Value = " "
$ws.Range("K24").Value = 0.21210000000000001
$ws.Range("L24").Value = 0.78059999999999996
$ws.Range("M24").Value = 0.0072
$ws.Range("N24").Formula = "=SUM(K24:M24)"

$ws.Range("J25").Value = "  "
$ws.Range("K25").Value = 0.4486
$ws.Range("L25").Value = 0.5514
$ws.Range("M25").Value = 0
$ws.Range("N25").Formula = "=SUM(K25:M25)"

$ws.Range("J26").Value = " "
$ws.Range("K26").Value = 0.67530000000000001
$ws.Range("L26").Value = 0.32469999999999999
$ws.Range("M26").Value = 0
$ws.Range("N26").Formula = "=SUM(K26:M26)"

$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Formula = "=SUM(K27:M27)"

# --- Clear marching-ants clipboard marquee left over from the copy operations ---
$excel.CutCopyMode = $false

# --- Update view state: scroll position and selection ---
$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
